$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Answer" column (G) placeholder option labels for rows 12-20.
$ws.Range("G12").Value = "Option – 1"
$ws.Range("G13").Value = "Option – 2"
$ws.Range("G14").Value = "Option – 12"
$ws.Range("G15").Value = "Option – 1"
$ws.Range("G16").Value = "Option – 1"
$ws.Range("G17").Value = "Option – 2"
$ws.Range("G18").Value = "Option – 4"
$ws.Range("G19").Value = "Option – 1"
$ws.Range("G20").Value = "Option – "

# Select the last-edited cell, matching the author's final cursor position.
$ws.Range("G20").Select()
